$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRADE")

# Fill column C (rows 3-37) with "NIL", matching the values already used
# elsewhere in the sheet (e.g. column B and the lower part of column C).
for ($r = 3; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value = "NIL"
}

# Update the view state: scrolled position and active selection.
$ws.Range("F31").Select()
$excel.ActiveWindow.ScrollRow = 41

$wb.Save()
